# Add a new paragraph after the "Email" paragraph containing a hyperlink
# whose display text is "https://github.com/JustinMKing".

$d = $word.ActiveDocument

# Anchor on the email address text (from the existing hyperlink) so the
# insertion point is found robustly rather than via a hard-coded index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*js.custom.carbon@gmail.com*") {
        $targetIndex = $i
    }
}

$emailPara = $d.Paragraphs.Item($targetIndex)

# Insert a brand new paragraph right after the Email paragraph. It inherits
# the paragraph mark's run formatting (rStyle "Emphasis"), matching the
# target paragraph's <w:pPr><w:rPr><w:rStyle w:val="Emphasis"/></w:rPr></w:pPr>.
$insertAt = $emailPara.Range
$insertAt.Collapse(0)
$insertAt.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range
$newRange.Collapse(1)

# Type the plain text first (this avoids leaving a stray empty run behind
# once it is converted into a hyperlink), then convert exactly that text
# range into a hyperlink run styled "Hyperlink".
$url = "https://github.com/JustinMKing"
$startPos = $newRange.Start
$newRange.InsertAfter($url)
$endPos = $startPos + $url.Length
$txtRange = $d.Range($startPos, $endPos)

$d.Hyperlinks.Add($txtRange, $url, [System.Reflection.Missing]::Value, `
    [System.Reflection.Missing]::Value, $url) | Out-Null
